$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Insert 4 new data rows right after the existing "period" row (row 16),
#    pushing the old row 17 down to row 21 and the signature rows (22/23)
#    down to (26/27). Excel keeps merged cells / references consistent.
# ---------------------------------------------------------------------------
$ws.Rows("17:20").Insert()

# ---------------------------------------------------------------------------
# 2) Clone the formatting (and content, temporarily) of row 16 into the four
#    freshly inserted rows so every cell picks up the correct style indexes
#    (borders, number format, fill) exactly like row 16 / the old row 17.
# ---------------------------------------------------------------------------
$ws.Range("B16:J16").Copy($ws.Range("B17:J17"))
$ws.Range("B16:J16").Copy($ws.Range("B18:J18"))
$ws.Range("B16:J16").Copy($ws.Range("B19:J19"))
$ws.Range("B16:J16").Copy($ws.Range("B20:J20"))

# ---------------------------------------------------------------------------
# 3) Existing worker (MIGUEL ANGEL SOLANO GOMEZ) now shows period 2401 first,
#    then the older 2312 period directly below (row 21 keeps the worker's
#    CC/name untouched from the copy above, only the period differs).
# ---------------------------------------------------------------------------
$ws.Range("E16").Value = "2401"
$ws.Range("E17").Value = "2312"

# ---------------------------------------------------------------------------
# 4) New worker MARELBY GONZALEZ RODRIGUEZ (CC 45504857) occupies the four
#    newly added periods 2505, 2504, 2503 and 2502, rows 18-21.
# ---------------------------------------------------------------------------
$ws.Range("C18").Value = "45504857"
$ws.Range("D18").Value = "MARELBY GONZALEZ RODRIGUEZ"
$ws.Range("E18").Value = "2505"
$ws.Range("F18").Value = 56940
$ws.Range("G18").Value = 1423500

$ws.Range("C19").Value = "45504857"
$ws.Range("D19").Value = "MARELBY GONZALEZ RODRIGUEZ"
$ws.Range("E19").Value = "2504"
$ws.Range("F19").Value = 56940
$ws.Range("G19").Value = 1423500

$ws.Range("C20").Value = "45504857"
$ws.Range("D20").Value = "MARELBY GONZALEZ RODRIGUEZ"
$ws.Range("E20").Value = "2503"
$ws.Range("F20").Value = 56940
$ws.Range("G20").Value = 1423500

$ws.Range("C21").Value = "45504857"
$ws.Range("D21").Value = "MARELBY GONZALEZ RODRIGUEZ"
$ws.Range("E21").Value = "2502"
$ws.Range("F21").Value = 56940
$ws.Range("G21").Value = 1423500

# ---------------------------------------------------------------------------
# 5) Update the summary block: total amount owed, worker count, period count.
# ---------------------------------------------------------------------------
$ws.Range("E11").Value = 320560
$ws.Range("C13").Value = 2
$ws.Range("F13").Value = 6
